# Generate Report for Handback
#
# This models a localization "handback" event: the zh-cn and de-de
# translation files for both source documents have come back in sync
# with en-US, so the status changes from "Ready for handoff" to
# "Handed back: in sync with en-US", the "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" columns on each
# language sheet get populated, and a couple of columns get widened so
# the new, longer values are readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status columns for both files (zh-cn / de-de) ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de status columns to fit the longer text.
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9920891c7827685675df55308e6ea0986f86b542/e2e/3537563b-77e7-4663-aa2e-957d039a8b2b.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "3537563b-77e7-4663-aa2e-957d039a8b2b.md")
$zhcn.Range("J2").Value = "3537563b-77e7-4663-aa2e-957d039a8b2b.85932e4da4ec79813dd07a19961c04b08e0771c4.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-13 23:03:13"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9920891c7827685675df55308e6ea0986f86b542/e2e/f5b1a7b0-c673-41c6-a019-b923d6cabc77.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "f5b1a7b0-c673-41c6-a019-b923d6cabc77.md")
$zhcn.Range("J3").Value = "f5b1a7b0-c673-41c6-a019-b923d6cabc77.5dbadb4022f65fb9cf506737aa36012f487f99e8.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-13 23:03:13"

# --- de-de sheet ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9920891c7827685675df55308e6ea0986f86b542/e2e/3537563b-77e7-4663-aa2e-957d039a8b2b.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "3537563b-77e7-4663-aa2e-957d039a8b2b.md")
$dede.Range("J2").Value = "3537563b-77e7-4663-aa2e-957d039a8b2b.85932e4da4ec79813dd07a19961c04b08e0771c4.de-de.xlf"
$dede.Range("K2").Value = "2016-08-13 23:03:22"

$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9920891c7827685675df55308e6ea0986f86b542/e2e/f5b1a7b0-c673-41c6-a019-b923d6cabc77.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "f5b1a7b0-c673-41c6-a019-b923d6cabc77.md")
$dede.Range("J3").Value = "f5b1a7b0-c673-41c6-a019-b923d6cabc77.5dbadb4022f65fb9cf506737aa36012f487f99e8.de-de.xlf"
$dede.Range("K3").Value = "2016-08-13 23:03:22"
